$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed per the repull of data.
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 6
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = -14
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -5
$ws.Range("F14").Value = 4
